# Normalize the "Recorded By" (column G) entries so that "System" (and,
# where no "System" entry is present, "admin@admin.com") is listed first
# among the comma-separated recorder names, matching the upstream sync.
#
# Examples:
#   "backup@backdoor.com, System, system"  -> "System, backup@backdoor.com, system"
#   "dnasr281@gmail.com, System"           -> "System, dnasr281@gmail.com"
#   "dnasr281@gmail.com, admin@admin.com"  -> "admin@admin.com, dnasr281@gmail.com"
# Rows where "System" (or "admin@admin.com") is already first are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Range("G$r")
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ", "
    if ($parts.Count -lt 2) { continue }

    $newVal = $val

    if ($parts[0] -eq "System") {
        # "System" already listed first - nothing to do.
        $newVal = $val
    }
    elseif ($parts[1] -eq "System") {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
        $newVal = $parts -join ", "
    }
    elseif ($parts[1] -eq "admin@admin.com" -and $parts[0] -ne "admin@admin.com") {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
        $newVal = $parts -join ", "
    }

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
